$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record at row 54 (pushes existing rows 54:119 down to 55:120,
# dimension grows from A1:R119 to A1:R120).
$ws.Rows("54:54").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44771
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112042
$ws.Range("G54").Value = "Locoto"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 150
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 22000
$ws.Range("M54").Value = 21000
$ws.Range("N54").Value = "`$/caja 20 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 1050
$ws.Range("Q54").Value = 20
$ws.Range("R54").Value = "Hortaliza"
